# Insert a new data row at row 25 (pushing the existing rows 25-120 down to
# 26-121) and populate the new row with a fresh "Granada" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(25).Insert()

$ws.Cells.Item(25, 1).Value  = 10
$ws.Cells.Item(25, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(25, 3).Value  = "La Araucanía"
$ws.Cells.Item(25, 4).Value  = 44701
$ws.Cells.Item(25, 5).Value  = 9
$ws.Cells.Item(25, 6).Value  = "Fruta"
$ws.Cells.Item(25, 7).Value  = 100104
$ws.Cells.Item(25, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(25, 9).Value  = 100104001
$ws.Cells.Item(25, 10).Value = "Granada"
$ws.Cells.Item(25, 11).Value = "Wonderfull"
$ws.Cells.Item(25, 12).Value = "Especial"
$ws.Cells.Item(25, 13).Value = 50
$ws.Cells.Item(25, 14).Value = 21000
$ws.Cells.Item(25, 15).Value = 21000
$ws.Cells.Item(25, 16).Value = 21000
$ws.Cells.Item(25, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(25, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 19).Value = 1400
$ws.Cells.Item(25, 20).Value = 15
